$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '26.193.50'
$ws.Range("E2").Value = '  -4.10%  '
$ws.Range("D3").Value = '1.659.28'
$ws.Range("E3").Value = '  -2.81%  '
Set-TextValue $ws.Range("D4") '1.006'
$ws.Range("E4").Value = '  +0.25%  '
Set-TextValue $ws.Range("D5") '217.66'
$ws.Range("E5").Value = '  -2.81%  '
Set-TextValue $ws.Range("D6") '0.5145'
$ws.Range("E6").Value = '  -3.37%  '
Set-TextValue $ws.Range("D8") '0.2581'
$ws.Range("E8").Value = '  -3.19%  '
Set-TextValue $ws.Range("D10") '20.00'
$ws.Range("E10").Value = '  -3.89%  '
Set-TextValue $ws.Range("D11") '0.07816'
$ws.Range("E11").Value = '  +2.28%  '
$ws.Range("D12").Value = '1.663.34'
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("E13").Value = '  -4.86%  '
$ws.Range("D14").Value = '1.887.47'
$ws.Range("E14").Value = '  -2.80%  '
Set-TextValue $ws.Range("D15") '0.5550'
$ws.Range("E15").Value = '  -4.72%  '
$ws.Range("D16").Value = '0.0₅8066'
$ws.Range("E16").Value = '  -1.42%  '
Set-TextValue $ws.Range("D17") '64.30'
$ws.Range("E17").Value = '  -4.92%  '
$ws.Range("D18").Value = '26.228.46'
$ws.Range("E18").Value = '  -4.12%  '
$ws.Range("E19").Value = '  +0.19%  '
Set-TextValue $ws.Range("D20") '211.07'
$ws.Range("E20").Value = '  -2.44%  '
Set-TextValue $ws.Range("D21") '4.429'
$ws.Range("E21").Value = '  -4.46%  '
Set-TextValue $ws.Range("D22") '10.09'
$ws.Range("E22").Value = '  -3.10%  '
Set-TextValue $ws.Range("D23") '6.059'
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("E24").Value = '  +0.19%  '
Set-TextValue $ws.Range("D25") '145.25'
$ws.Range("E25").Value = '  +0.99%  '
Set-TextValue $ws.Range("D26") '1.755'
$ws.Range("E26").Value = '  +2.72%  '
$ws.Range("E27").Value = '  -2.70%  '
Set-TextValue $ws.Range("D28") '6.995'
$ws.Range("E28").Value = '  -3.19%  '
$ws.Range("E29").Value = '  -2.48%  '
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("E31").Value = '  -2.73%  '
Set-TextValue $ws.Range("D32") '3.359'
$ws.Range("E32").Value = '  -3.37%  '
$ws.Range("E33").Value = '  -5.23%  '
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("E35").Value = '  -3.66%  '
Set-TextValue $ws.Range("D36") '0.9312'
$ws.Range("E36").Value = '  -1.98%  '
Set-TextValue $ws.Range("D37") '2.374'
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("D38").Value = '1.172.65'
$ws.Range("E38").Value = '  +12.33%  '
Set-TextValue $ws.Range("D39") '0.5705'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  -2.64%  '
$ws.Range("E41").Value = '  +0.19%  '
Set-TextValue $ws.Range("D42") '0.8401'
$ws.Range("E42").Value = '  -0.29%  '
Set-TextValue $ws.Range("D43") '5.682'
$ws.Range("E43").Value = '  -2.18%  '
Set-TextValue $ws.Range("D44") '100.68'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").Value = '1.797.76'
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("E46").Value = '  +4.58%  '
Set-TextValue $ws.Range("D47") '0.4541'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("E48").Value = '  -3.31%  '
$ws.Range("E49").Value = '  -0.05%  '
Set-TextValue $ws.Range("D50") '7.893'
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("E51").Value = '  -3.23%  '
